# Commit: "case as member of doc"
#
# The template referenced `case` directly (e.g. `{{case.footer}}`,
# `{{case.cause_number}}`). The fix re-parents it as a member of `doc`
# (e.g. `{{doc.case.footer}}`, `{{doc.case.cause_number}}`).
#
# Both occurrences live in the page footer. The footer run carrying
# "case.footer" is rendered in small caps (w:caps), so Find must use
# MatchCase = $false to locate the lower-case stored text even though
# it is visually all-caps; the replacement text is supplied in the
# correct (lower) case and is inserted verbatim regardless of MatchCase.

$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {
    foreach ($footer in $sec.Footers) {
        $rng = $footer.Range

        $rng.Find.Execute(
            "{{case.footer}}", $false, $false, $false, $false, $false,
            $true, 1, $false, "{{doc.case.footer}}", 2) | Out-Null

        $rng.Find.Execute(
            "{{case.cause_number}}", $false, $false, $false, $false, $false,
            $true, 1, $false, "{{doc.case.cause_number}}", 2) | Out-Null
    }
}

Write-Host "Done."
